# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# worksheets to match the latest scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new F-column value. Both sheets share identical data (aside from
# the already-differing F8 seed value), so the same set of updates applies
# to each sheet.
$updates = @{
    3  = 82
    6  = 516
    7  = 45
    8  = 1961
    10 = 86
    11 = 4135
    13 = 269
    14 = 95
    15 = 86
    17 = 52
    18 = 2849
    20 = 399
    23 = 60
    25 = 53
    28 = 41
    29 = 186
    30 = 271
    31 = 1618
    32 = 228
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
